# Add "Wins" / "Losses" / "Ties" season-record columns (AD:AF) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) — same bold/bordered/centered style as the existing headers.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AD1:AF1").Font.Bold = $true
$ws.Range("AD1:AF1").HorizontalAlignment = -4108
$ws.Range("AD1:AF1").VerticalAlignment = -4160
$ws.Range("AD1:AF1").Borders.LineStyle = 1

# Data rows (2-46) — every player row gets the same season record.
$lastRow = 46
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 61
    $ws.Cells.Item($r, 31).Value = 101
    $ws.Cells.Item($r, 32).Value = 0
}
